# Apply updates to the "Eventos" sheet reflecting the refreshed daily events export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (vehicleId), E (Unidad) and F (driverId) hold numeric-looking
# identifiers that must stay stored as text (as in the source export), so
# force a Text number format on those cells before writing the values.
$ws.Range("D2:F4").NumberFormat = "@"

# --- Row 2 ---
$ws.Range("A2").Value = "281474991265569-1739490786629"
$ws.Range("B2").Value = "Harsh Brake"
$ws.Range("C2").Value = "2025-02-13T17:53:06.629"
$ws.Range("D2").Value = "281474991265569"
$ws.Range("E2").Value = "103"
$ws.Range("F2").Value = "52215661"
$ws.Range("G2").Value = "CARLOS ALBERTO JIMENEZ"
$ws.Range("H2").Value = 20.747010889
$ws.Range("I2").Value = -103.398126309
$ws.Range("J2").Value = 0.7530547976493835
$ws.Range("K2").Value = "No video URL"

# --- Row 3 ---
$ws.Range("A3").Value = "281474991265569-1739486708088"
$ws.Range("B3").Value = "Harsh Brake"
$ws.Range("C3").Value = "2025-02-13T16:45:08.088"
$ws.Range("D3").Value = "281474991265569"
$ws.Range("E3").Value = "103"
$ws.Range("F3").Value = "52215661"
$ws.Range("G3").Value = "CARLOS ALBERTO JIMENEZ"
$ws.Range("H3").Value = 20.697803669
$ws.Range("I3").Value = -103.386574499
$ws.Range("J3").Value = 0.8436447978019714
$ws.Range("K3").Value = "No video URL"

# --- Row 4 ---
$ws.Range("A4").Value = "281474991109374-1739468687126"
$ws.Range("B4").Value = "No Seat Belt"
$ws.Range("C4").Value = "2025-02-13T11:44:47.126"
$ws.Range("D4").Value = "281474991109374"
$ws.Range("E4").Value = "102"
$ws.Range("F4").Value = "51848506"
$ws.Range("G4").Value = "ARMANDO MUÑOZ"
$ws.Range("H4").Value = 20.293003759
$ws.Range("I4").Value = -102.619195559
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991109374/1739468684626/xJ0Yy0MOUz-camera-video-segment-driver-1739468687126.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHMRWCAAS%2F20250214%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250214T151445Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEAQaCXVzLXdlc3QtMiJHMEUCIQDVg77hQIY9CiYyCDZb%2FCGqL78BfW6ZPneQUULMOofxawIgM5Zt8BAsGXb%2BaDNE6jMg8aPFy2U05i7DAcoXS1yx17Yq3QMILBAEGgw3ODEyMDQ5NDIyNDQiDGDfN71NTKfiRfHueSq6A1ukKCo2nB5syXluJjCSoBe0iqheLrvqrdQCTM0mRFYHySa1Bv9j0D5Ukl1PFc%2F9%2FF7dY5nTXN8d8NM3Z2RdmZx5Af1QE2ES8ZihKdmPBdc4ELUSCsTbFJDaqC7dMjQ7WGXRN%2F5zgrmq5YZJyWVuZBJDQRONg62ReeCjIqLAy%2BZsGfCZJYmxCkK3vphISzB1HJD9snnHOyTfeFHDvnKqMSpondsqkkxGkVMrfAJiA227Pphch%2BzzAuSz%2FKwG7nt52PDa4l60JvsrHFQonJOE5E%2B8yQnsAtdUDetRCvR0sQGb7JCFuGB7APjNMTt0dF0JK9r%2FsBCR9BEoDQwsqvQ5VaJIUo0bcU%2BUkm9vZwU6UH8Rn%2BFMcOJBl%2Bf0eHWIvCUAyTeAvxdFyTDvIVmnFazHQE%2BZt1lumm0aqWr6GL10X9yovBLNxoCFyaHugv5bAL0W46IcoFNqdzjunZTevxAdyYiUQ9r%2FysRvtJs3Fc6ULO6ABeQ%2BN%2Fn00mVBdPSRQGuhEIkS1hrxxE3UBDeqqIIWrniu8NnudTLp7%2FByQnnP58yjyc5EKl4LimI5V%2FBMcMM09JW9FnKe1QvT71sw4s28vQY6pQEyiWEJM6aIqQyG3rgkVzpdz2ZQQtvt6v2%2BKAl4DCIF7nOHeMgWuvMshoJDB%2F1cT8HDQhVCp8LbYIDR2zeaBzKar1iE3%2Bl0DR1FdNmDyyW83Y%2F%2BOVxMmZbuvdfzG3t8WAiNAA2mmtjgR%2B6NKqCt9J8%2Fo0l0nAW32Yoi8DCLQySNhuCYhCObOI5BPBoATpnn9%2FJuLugFDnTK7XTlfb%2FWihEm2C3qemE%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2014%20Feb%202025%2023%3A14%3A45%20GMT&X-Amz-Signature=32f7283401873bcb1a91694f8094beb564dc1b6107365f2b05ff503ca8824556"
$ws.Range("L4").Value = "No video URL"
